# SoIB_summaries.xlsx update after re-running resolve / classify+summarise
# steps following changes to the mapping file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# "Range Status" sheet: species counts now all zero, and the
# percentage column (C) is no longer populated for any row.
# ---------------------------------------------------------------
$wsRange = $wb.Worksheets.Item("Range Status")
for ($r = 2; $r -le 7; $r++) {
    $wsRange.Cells.Item($r, 2).Value = 0
    $wsRange.Cells.Item($r, 3).ClearContents()
}

# ---------------------------------------------------------------
# "Species qualification" sheet: Range Analysis species count
# dropped to 0.
# ---------------------------------------------------------------
$wsQual = $wb.Worksheets.Item("Species qualification")
$wsQual.Range("B5").Value = 0

# ---------------------------------------------------------------
# "High Priority break-up" sheet: the breakdown changed -
# "Trend Different" became "Trend New" (with a new percentage),
# the old "Range" row data was replaced with what used to be the
# "IUCN" row's label but new numbers, and the trailing IUCN row
# was removed entirely.
# ---------------------------------------------------------------
$wsBreak = $wb.Worksheets.Item("High Priority break-up")

$wsBreak.Range("A2").Value = "Trend New"
$wsBreak.Range("E2").Value = 4.2

$wsBreak.Range("A3").Value = "IUCN"
$wsBreak.Range("B3").Value = 23
$wsBreak.Range("C3").Value = 95.8
$wsBreak.Range("D3").Value = 23
$wsBreak.Range("E3").Value = 95.8

$wsBreak.Range("A4:E4").ClearContents()
